$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A10 currently holds the text "26/5/2015" - convert it to a real date value
# matching the style used by the other date cells above (A4:A9), i.e. number
# format applied via style index 1 (numFmtId 14).
$ws.Range("A10").Value = 42150
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)

# Add the new row of data (row 11)
$ws.Range("A11").Value = 42151
$ws.Range("A9").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("B11").Value = "Learn SASS and use SASS to create CSS for the demo project"

# Update the active selection to match the new end-of-data location
$ws.Range("B12").Select()
